# Update "想去人数" (interested-people count) figures in column F for the
# 展览 (Exhibition) and 全部类型 (All Types) sheets to the freshly scraped
# totals from the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3538
$ws1.Range("F3").Value = 747
$ws1.Range("F4").Value = 148
$ws1.Range("F5").Value = 7033
$ws1.Range("F6").Value = 3398
$ws1.Range("F7").Value = 59
$ws1.Range("F8").Value = 156
$ws1.Range("F12").Value = 47
$ws1.Range("F13").Value = 25
$ws1.Range("F14").Value = 181
$ws1.Range("F15").Value = 595
$ws1.Range("F16").Value = 40
$ws1.Range("F17").Value = 46

# --- Sheet "全部类型" -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3538
$ws4.Range("F4").Value = 747
$ws4.Range("F5").Value = 149
$ws4.Range("F6").Value = 7033
$ws4.Range("F7").Value = 3398
$ws4.Range("F8").Value = 59
$ws4.Range("F9").Value = 156
$ws4.Range("F13").Value = 47
$ws4.Range("F14").Value = 25
$ws4.Range("F15").Value = 181
$ws4.Range("F16").Value = 595
$ws4.Range("F17").Value = 40
$ws4.Range("F18").Value = 46
